$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (Player, Position, Team)
$data = @(
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Patrick Williams", "PF", "Chicago Bulls"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Naji Marshall", "SG,SF", "Dallas Mavericks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Walker Kessler", "C", "Utah Jazz")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
